$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 477.2353
$ws.Range("I19").Value = 226.375
$ws.Range("K19").Value = 226.375
$ws.Range("M19").Value = -51.375
$ws.Range("H74").Value = 18414.143
$ws.Range("I74").Value = 18414.143
$ws.Range("K74").Value = 18414.143
$ws.Range("M74").Value = -17478.143
$ws.Range("H77").Value = 18414.143
$ws.Range("I77").Value = 18414.143
$ws.Range("K77").Value = 92070.715
$ws.Range("M77").Value = -87390.715
$ws.Range("H132").Value = 4055.9768
$ws.Range("I132").Value = 4435.6216
$ws.Range("K132").Value = 13306.8648
$ws.Range("M132").Value = -10776.8648
$ws.Range("H137").Value = 2499.7576
$ws.Range("I137").Value = 1931.2222
$ws.Range("K137").Value = 5793.6666
$ws.Range("M137").Value = -3243.6666
$ws.Range("H138").Value = 266593.5
$ws.Range("I138").Value = 3609.8333
$ws.Range("J138").Value = 387970.6
$ws.Range("K138").Value = 10829.4999
$ws.Range("L138").Value = 1163911.8
$ws.Range("M138").Value = -5689.499899999999
$ws.Range("N138").Value = -1174191.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1943.82
$ws.Range("I32").Value = 1718.1837
$ws.Range("J32").Value = 13000
$ws.Range("K32").Value = 1718.1837
$ws.Range("L32").Value = 13000
$ws.Range("M32").Value = -1431.1837
$ws.Range("N32").Value = -13574
$ws.Range("H61").Value = 3804.7144
$ws.Range("I61").Value = 2230.9285
$ws.Range("J61").Value = 6952.2856
$ws.Range("K61").Value = 2230.9285
$ws.Range("L61").Value = 6952.2856
$ws.Range("M61").Value = -2018.9285
$ws.Range("N61").Value = -7376.2856
$ws.Range("H122").Value = 4606.5483
$ws.Range("I122").Value = 3921.5
$ws.Range("J122").Value = 8168.8
$ws.Range("K122").Value = 11764.5
$ws.Range("L122").Value = 24506.4
$ws.Range("M122").Value = -9314.5
$ws.Range("N122").Value = -29406.4
$ws.Range("H136").Value = 3804.7144
$ws.Range("I136").Value = 2230.9285
$ws.Range("J136").Value = 6952.2856
$ws.Range("K136").Value = 6692.7855
$ws.Range("L136").Value = 20856.8568
$ws.Range("M136").Value = -4142.7855
$ws.Range("N136").Value = -25956.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 64814.668
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H86").Value = 5663.6665
$ws.Range("I86").Value = 5412.3335
$ws.Range("J86").Value = 6166.3335
$ws.Range("K86").Value = 5412.3335
$ws.Range("L86").Value = 6166.3335
$ws.Range("M86").Value = -4289.3335
$ws.Range("N86").Value = -8412.333500000001
$ws.Range("H89").Value = 5663.6665
$ws.Range("I89").Value = 5412.3335
$ws.Range("J89").Value = 6166.3335
$ws.Range("K89").Value = 27061.6675
$ws.Range("L89").Value = 30831.6675
$ws.Range("M89").Value = -21445.6675
$ws.Range("N89").Value = -42063.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3176.6582
$ws.Range("I31").Value = 2652.678
$ws.Range("K31").Value = 2652.678
$ws.Range("M31").Value = -2357.678
$ws.Range("H34").Value = 3176.6582
$ws.Range("I34").Value = 2652.678
$ws.Range("K34").Value = 2652.678
$ws.Range("M34").Value = -2450.678
$ws.Range("H58").Value = 2916.56
$ws.Range("I58").Value = 1803.5
$ws.Range("J58").Value = 4895.3335
$ws.Range("K58").Value = 1803.5
$ws.Range("L58").Value = 4895.3335
$ws.Range("M58").Value = -1600.5
$ws.Range("N58").Value = -5301.3335
$ws.Range("H107").Value = 779.7917
$ws.Range("I107").Value = 811.3158
$ws.Range("J107").Value = 660
$ws.Range("K107").Value = 811.3158
$ws.Range("L107").Value = 660
$ws.Range("M107").Value = 1108.6842
$ws.Range("N107").Value = -4500
$ws.Range("H122").Value = 3444.842
$ws.Range("I122").Value = 3261.8823
$ws.Range("K122").Value = 9785.6469
$ws.Range("M122").Value = -7335.6469
$ws.Range("H132").Value = 3614.4167
$ws.Range("I132").Value = 3860.3333
$ws.Range("J132").Value = 3122.5833
$ws.Range("K132").Value = 11580.9999
$ws.Range("L132").Value = 9367.749899999999
$ws.Range("M132").Value = -9050.999899999999
$ws.Range("N132").Value = -14427.7499
$ws.Range("H134").Value = 2036.2084
$ws.Range("I134").Value = 1609.921
$ws.Range("J134").Value = 3656.1
$ws.Range("K134").Value = 4829.763
$ws.Range("L134").Value = 10968.3
$ws.Range("M134").Value = -2294.763
$ws.Range("N134").Value = -16038.3
$ws.Range("H136").Value = 2916.56
$ws.Range("I136").Value = 1803.5
$ws.Range("J136").Value = 4895.3335
$ws.Range("K136").Value = 5410.5
$ws.Range("L136").Value = 14686.0005
$ws.Range("M136").Value = -2860.5
$ws.Range("N136").Value = -19786.0005
$ws.Range("H141").Value = 294849.8
$ws.Range("J141").Value = 294849.8
$ws.Range("L141").Value = 294849.8
$ws.Range("N141").Value = -305209.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7460.75
$ws.Range("I56").Value = 7460.75
$ws.Range("K56").Value = 7460.75
$ws.Range("M56").Value = -6930.75
$ws.Range("H75").Value = 496.75
$ws.Range("J75").Value = 553.4286
$ws.Range("L75").Value = 1660.2858
$ws.Range("N75").Value = -3656.2858
$ws.Range("H78").Value = 496.75
$ws.Range("J78").Value = 553.4286
$ws.Range("L78").Value = 4980.8574
$ws.Range("N78").Value = -14964.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6040.6523
$ws.Range("I102").Value = 1751.7222
$ws.Range("K102").Value = 1751.7222
$ws.Range("M102").Value = -129.7221999999999
$ws.Range("H103").Value = 90000
$ws.Range("J103").Value = 90000
$ws.Range("L103").Value = 90000
$ws.Range("N103").Value = -92344
$ws.Range("H135").Value = 64284.57
$ws.Range("I135").Value = 30998
$ws.Range("J135").Value = 69832.336
$ws.Range("K135").Value = 30998
$ws.Range("L135").Value = 69832.336
$ws.Range("M135").Value = -25928
$ws.Range("N135").Value = -79972.336
$ws.Range("H138").Value = 102999
$ws.Range("J138").Value = 102999
$ws.Range("L138").Value = 102999
$ws.Range("N138").Value = -113279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1728
$ws.Range("H27").Value = 1728
$ws.Range("H132").Value = 4190.727
$ws.Range("I132").Value = 3512.375
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 10537.125
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -8007.125
$ws.Range("N132").Value = -23058.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 985.4737
$ws.Range("I113").Value = 1240.1111
$ws.Range("J113").Value = 756.3
$ws.Range("K113").Value = 3720.3333
$ws.Range("L113").Value = 2268.9
$ws.Range("M113").Value = -1550.3333
$ws.Range("N113").Value = -6608.9
$ws.Range("H136").Value = 35715852
$ws.Range("I136").Value = 43479776
$ws.Range("J136").Value = 1799.4
$ws.Range("K136").Value = 130439328
$ws.Range("L136").Value = 5398.200000000001
$ws.Range("M136").Value = -130436778
$ws.Range("N136").Value = -10498.2
$ws.Range("H141").Value = 69544.37
$ws.Range("J141").Value = 69544.37
$ws.Range("L141").Value = 69544.37
$ws.Range("N141").Value = -79904.37
